# Insert two new rows (256 and 257) for the week of date 44615 (2022-02-23),
# shifting all following rows down by two. This mirrors Excel's "Insert"
# row behavior (entire row insert), which pushes the previous rows 256..273
# down to 258..275 and updates the sheet dimension automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A256:A257").EntireRow.Insert()

# New row 256: Especial quality entry for 2022-02-23
$ws.Cells.Item(256, 1).Value = 3
$ws.Cells.Item(256, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(256, 3).Value = "Coquimbo"
$ws.Cells.Item(256, 4).Value = 44615
$ws.Cells.Item(256, 5).Value = 5
$ws.Cells.Item(256, 6).Value = "Fruta"
$ws.Cells.Item(256, 7).Value = 100101
$ws.Cells.Item(256, 8).Value = "Berries"
$ws.Cells.Item(256, 9).Value = 100112025
$ws.Cells.Item(256, 10).Value = "Frutilla"
$ws.Cells.Item(256, 11).Value = "Sin especificar"
$ws.Cells.Item(256, 12).Value = "Especial"
$ws.Cells.Item(256, 13).Value = 56
$ws.Cells.Item(256, 14).Value = 6000
$ws.Cells.Item(256, 15).Value = 6000
$ws.Cells.Item(256, 16).Value = 6000
$ws.Cells.Item(256, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(256, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(256, 19).Value = 857
$ws.Cells.Item(256, 20).Value = 7

# New row 257: Segunda quality entry for 2022-02-23
$ws.Cells.Item(257, 1).Value = 3
$ws.Cells.Item(257, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(257, 3).Value = "Coquimbo"
$ws.Cells.Item(257, 4).Value = 44615
$ws.Cells.Item(257, 5).Value = 5
$ws.Cells.Item(257, 6).Value = "Fruta"
$ws.Cells.Item(257, 7).Value = 100101
$ws.Cells.Item(257, 8).Value = "Berries"
$ws.Cells.Item(257, 9).Value = 100112025
$ws.Cells.Item(257, 10).Value = "Frutilla"
$ws.Cells.Item(257, 11).Value = "Sin especificar"
$ws.Cells.Item(257, 12).Value = "Segunda"
$ws.Cells.Item(257, 13).Value = 48
$ws.Cells.Item(257, 14).Value = 4000
$ws.Cells.Item(257, 15).Value = 4000
$ws.Cells.Item(257, 16).Value = 4000
$ws.Cells.Item(257, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(257, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(257, 19).Value = 571
$ws.Cells.Item(257, 20).Value = 7

Write-Host "Inserted rows 256-257 with new Frutilla price entries for 2022-02-23"
